# Literature Survey document edit:
# After the existing "... framework for Deep Learning" (Caffe) bullet point,
# insert a page-break section introducing a new top-level topic, "The Weather
# Classification problem", with its "Description" (a quoted excerpt),
# "Applications" and "Approaches" sub-bullets, followed by another page break.

$d = $word.ActiveDocument

# Anchor on the paragraph that ends the existing bullet list ("... framework
# for Deep Learning") using Find, so this does not depend on absolute
# paragraph indices.
$anchorRange = $d.Content
$anchorRange.Find.Execute("framework for Deep Learning", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $anchorRange.Paragraphs.First

# Insert a fresh paragraph right after it; $anchorPara.Index stays valid
# (it is not renumbered by the insert), so the new blank paragraph is at
# index ($anchorPara.Index + 1).
$anchorPara.Range.InsertParagraphAfter()
$newIndex = $anchorPara.Index + 1
$target = $d.Paragraphs.Item($newIndex)

# Replace that blank paragraph's content/formatting wholesale with the exact
# OOXML for all of the new material (this also clears any pPr/numPr the
# blank paragraph might have inherited from its neighbour).
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:br/></w:r></w:p><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>The Weather Classification problem</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Description</w:t></w:r></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240" w:line="300" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr><w:r><w:t>“</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Vision based driver assistance systems (DAS) are currently designed to perform under good-natured weather conditions. Unfortunately, limited visibility often occurs in daily life (e.g. heavy rain or fog). As this strongly affects the accuracy or even the general function of vision systems, the actual weather condition is </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>a valuable</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> information for assistance systems. Based on the results of weather classification, specialized approaches for each class can be invoked to improve cognition. This will form a key factor to expand the application of DAS from selected environmental conditions to an overall approach.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">” </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Applications</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Approaches</w:t></w:r></w:p><w:p/><w:p><w:r><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)

Write-Output "Paragraphs after insert: $($d.Paragraphs.Count)"
